$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 771
$ws1.Range("F6").Value = 129
$ws1.Range("F8").Value = 136
$ws1.Range("F10").Value = 441
$ws1.Range("F11").Value = 500
$ws1.Range("F13").Value = 11535
$ws1.Range("F14").Value = 5388

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 771
$ws4.Range("F8").Value = 129
$ws4.Range("F10").Value = 136
$ws4.Range("F12").Value = 441
$ws4.Range("F13").Value = 500
$ws4.Range("F15").Value = 11535
$ws4.Range("F17").Value = 5388
